$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell values per the diff
$ws.Range("B2").Value = 121
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 43

# Remove the old row 4 entirely (was A4=1, B4=43)
$ws.Rows.Item(4).Delete()
